# Update betting odds values in Sheet1 per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("N3").Value = 1.5
$ws.Range("O3").Value = 2.5
$ws.Range("T4").Value = 8.5
$ws.Range("Y4").Value = 41
$ws.Range("AF4").Value = 11
$ws.Range("J6").Value = 1.06
$ws.Range("K6").Value = 10
$ws.Range("L6").Value = 1.33
$ws.Range("M6").Value = 3.4
$ws.Range("N6").Value = 2
$ws.Range("O6").Value = 1.8
$ws.Range("J9").Value = 1.03
$ws.Range("L9").Value = 1.25
$ws.Range("G10").Value = 1.85
$ws.Range("I10").Value = 4.1
$ws.Range("J10").Value = 1.03
$ws.Range("L10").Value = 1.19
$ws.Range("R10").Value = 1.7
$ws.Range("S10").Value = 2.05
$ws.Range("T10").Value = 8
$ws.Range("W10").Value = 15
$ws.Range("J13").Value = 1.05
$ws.Range("K13").Value = 9
$ws.Range("L13").Value = 1.37
$ws.Range("M13").Value = 2.75
$ws.Range("N13").Value = 2.3
$ws.Range("O13").Value = 1.6
$ws.Range("L14").Value = 1.36
$ws.Range("M14").Value = 2.9
$ws.Range("R14").Value = 2.37
$ws.Range("S14").Value = 1.52
$ws.Range("H16").Value = 3.65
$ws.Range("I16").Value = 4.5
$ws.Range("L16").Value = 1.29
$ws.Range("M16").Value = 3.3
$ws.Range("R16").Value = 1.86
$ws.Range("S16").Value = 1.84
$ws.Range("T16").Value = 5.7
$ws.Range("U16").Value = 6.3
$ws.Range("W16").Value = 10
$ws.Range("AC16").Value = 60
$ws.Range("AE16").Value = 10
$ws.Range("AF16").Value = 20
$ws.Range("AH16").Value = 55
$ws.Range("AI16").Value = 35
$ws.Range("G19").Value = 2.22
$ws.Range("H19").Value = 3.4
$ws.Range("I19").Value = 2.87
$ws.Range("L19").Value = 1.21
$ws.Range("O19").Value = 2
$ws.Range("R19").Value = 1.53
$ws.Range("S19").Value = 2.18
$ws.Range("T19").Value = 9.75
$ws.Range("U19").Value = 12.5
$ws.Range("V19").Value = 8.75
$ws.Range("W19").Value = 23
$ws.Range("X19").Value = 16.5
$ws.Range("AA19").Value = 6.8
$ws.Range("AB19").Value = 12
$ws.Range("AE19").Value = 11.25
$ws.Range("AF19").Value = 16.5
$ws.Range("AG19").Value = 10.25
$ws.Range("AH19").Value = 35
$ws.Range("AI19").Value = 22
$ws.Range("AJ19").Value = 26
$ws.Range("G21").Value = 2.52
$ws.Range("H21").Value = 3
$ws.Range("I21").Value = 2.77
$ws.Range("L21").Value = 1.33
$ws.Range("M21").Value = 2.8
$ws.Range("N21").Value = 1.98
$ws.Range("O21").Value = 1.65
$ws.Range("P21").Value = 1.39
$ws.Range("Q21").Value = 2.55
$ws.Range("S21").Value = 1.9
$ws.Range("U21").Value = 12.5
$ws.Range("W21").Value = 28
$ws.Range("X21").Value = 22
$ws.Range("Y21").Value = 32
$ws.Range("Z21").Value = 8.5
$ws.Range("AA21").Value = 5.8
$ws.Range("AB21").Value = 13.5
$ws.Range("AG21").Value = 10
$ws.Range("AH21").Value = 35
$ws.Range("I22").Value = 2.6
$ws.Range("Z23").Value = 10.75
$ws.Range("K24").Value = 17
$ws.Range("R24").Value = 1.47
$ws.Range("O26").Value = 1.77
$ws.Range("AA27").Value = 5.9
$ws.Range("J33").Value = 1.08
$ws.Range("K33").Value = 8
